# Applies the BOM update for Chit_Chat/Probability_Gate_2022_BOM.xlsx
#  - Bumps the document version note
#  - Renames CD4093/CD4013 references to HEF4093/HEF4013 and adds remarks
#  - Fixes a resistor value typo (20k -> 22k)
#  - Replaces the IDC shrouded box header part with plain pin headers, with a remark
#  - Moves the selection to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probability_Gate_2022")

# Header note: document version bump
$ws.Range("B6").Value = "Document Version 03/02/2023"

# R17, R18 value correction
$ws.Range("C23").Value = "22k"

# U8, U9 : CD4093 -> HEF4093, clear Tayda part number, add remark
$ws.Range("C28").Value = "HEF4093"
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = "Must be HEF4093. The CD… version has different thresholds for schmitt trigger"

# U10 : CD4013 -> HEF4013, clear Tayda part number, add remark
$ws.Range("C29").Value = "HEF4013"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = "Can alsobe CD4013 (not tested)"

# U11 : IDC Shrouded Box Header -> POWER (plain pin headers), add remark
$ws.Range("C31").Value = "POWER"
$ws.Range("D31").Value = "Male Pin headers 2x05_P2.54mm_Vertical"
$ws.Range("F31").Value = "Shrouded IDC doesn't really fit"

# Move active selection
$ws.Range("C6").Select()
